$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "The modern world relies heavily on air travel" -> "... on-air travel"
#    In the source XML "on air" sits in its own run, wrapped in
#    <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
#    markers. Replacing a span that crosses into the neighbouring runs
#    makes Word regenerate that part of the paragraph, which drops the
#    now-stale proofing marks along with it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("heavily on air travel", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "heavily on-air travel", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Drop the closing "---" separator, the blank paragraph after it, and
#    the final "This report provides a concise overview ..." paragraph
#    (that summary now lives in a separate results/presentation file).
#    The empty paragraph right before "---" is left in place and becomes
#    the document's last paragraph.
# ---------------------------------------------------------------------
$sepRange = $d.Content
$found = $sepRange.Find.Execute("---", $true, $false, $false, $false, `
                                 $false, $true, 1, $false, "", 0)
if ($found -and $sepRange.Find.Found) {
    $startPara = $sepRange.Paragraphs(1)
    $lastPara = $d.Paragraphs.Last
    $killRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
    $killRange.Delete()
}
